$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The old column H ("date_of_creation") is dropped entirely in the new
# schema, so clear it out before laying down the new 7-column table.
# ---------------------------------------------------------------------------
$ws.Range("H1:H3").Clear() | Out-Null

# ---------------------------------------------------------------------------
# Header row (row 1) — relabel / reorder the schema columns (7 cols, was 8)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "creation_date"
$ws.Range("B1").Value = "currency"
$ws.Range("C1").Value = "settlement_amount"
$ws.Range("D1").Value = "buy_sell"
$ws.Range("E1").Value = "isin"
$ws.Range("F1").Value = "settlement_date"
$ws.Range("G1").Value = "SSI"

# ---------------------------------------------------------------------------
# Data rows
# ---------------------------------------------------------------------------

# Row 2
$ws.Range("A2").Value = 45252
$ws.Range("B2").Value = "USD"
$ws.Range("C2").Value = 11849200
$ws.Range("D2").Value = "buy"
$ws.Range("E2").Value = "US6937512345"
$ws.Range("F2").Value = 45983
$ws.Range("G2").Value = "Stream A/c 19200 Our Beneficiary Swift Code BKBKBK32 Our Agent BIC Code BKBKBKBK"

# Row 3
$ws.Range("A3").Value = 45220
$ws.Range("B3").Value = "USD"
$ws.Range("C3").Value = 19777430.56
$ws.Range("D3").Value = "sell"
$ws.Range("E3").Value = "US9127123213"
$ws.Range("F3").Value = 45951
$ws.Range("G3").Value = "PSET FFFF33"

# Row 4 (new)
$ws.Range("A4").Value = 45930
$ws.Range("B4").Value = "USD"
$ws.Range("C4").Value = 29851455.46
$ws.Range("D4").Value = "buy"
$ws.Range("E4").Value = "US912812313"
$ws.Range("F4").Value = 45931
$ws.Range("G4").Value = "Our Settlement Instructions BANK OF NEW YORK, NEW YORK (BDS) FXF  Your Settlement Instructions CITBNK N.A. 0854CS104"

# Row 5 (new)
$ws.Range("A5").Value = 45930
$ws.Range("B5").Value = "EUR"
$ws.Range("C5").Value = 12312960
$ws.Range("D5").Value = "buy"
$ws.Range("E5").Value = "FR00123123"
$ws.Range("F5").Value = 45932
$ws.Range("G5").Value = "GOLDMAN SACHS`n/ER/1111`nCLIENT`n/ER/2222`n"

# Row 6 (new) — no buy_sell value for this row
$ws.Range("A6").Value = 45978
$ws.Range("B6").Value = "USD"
$ws.Range("C6").Value = 33862165.87
$ws.Range("D6").Value = $null
$ws.Range("E6").Value = "US9127111117"
$ws.Range("F6").Value = 45978
$ws.Range("G6").Value = "0213123131089 CINK NYC`nCUST`nA/C 11111`nXXX12323`nOur Settlement:`nABA # 03123123118 BK OF NYC/1231313, CIOUP GLOBAL MARKETS INC AS SETTLEMENT AGENT FOR 1231"

# ---------------------------------------------------------------------------
# Number formats
# ---------------------------------------------------------------------------

# Date columns (A and F) — reuse the workbook's existing date style (already
# applied to A2/A3) by copying formats across, rather than minting new ones.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A4:A6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A2").Copy() | Out-Null
$ws.Range("F2:F6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Settlement amount column: #,##0.00 on most rows, accounting/comma style on row 3
$ws.Range("C2").NumberFormat = "#,##0.00"
$ws.Range("C4:C6").NumberFormat = "#,##0.00"
$ws.Range("C3").NumberFormat = "_(* #,##0.00_);_(* \(#,##0.00\);_(* ""-""??_);_(@_)"

# ---------------------------------------------------------------------------
# Alignment quirks on the SSI column for the two newest multi-line rows
# ---------------------------------------------------------------------------
$ws.Range("G5").HorizontalAlignment = -4131
$ws.Range("G6").WrapText = $false

# ---------------------------------------------------------------------------
# Column widths (approximate autofit results for the new content)
# ---------------------------------------------------------------------------
$ws.Range("A1").ColumnWidth = 16.28515625
$ws.Range("F1").ColumnWidth = 16

# ---------------------------------------------------------------------------
# View / selection
# ---------------------------------------------------------------------------
$ws.Range("A6").Select()

Write-Output "edit applied"
